$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (Year / 2019 / 2018 / 2017) is left untouched - it already holds
# the correct values and its bold/bordered style (s="1").

# New data for columns B..J. Each entry is the column's 4 values
# (row1 header, row2, row3, row4). All values are written as literal TEXT
# (not numbers) to match the source table, which stores every cell -
# including years and numeric-looking figures - as shared strings.
$columns = @(
    @("B", @("Operating profit", "375", "499", "598")),
    @("C", @("Profit before taxes", "460", "310", "320")),
    @("D", @("Taxes", "-66", "12", "250")),
    @("E", @("Net income for the year", "394", "322", "570")),
    @("F", @("Total assets", "284,305", "258,548", "251,998")),
    @("G", @("Equity", "10,576", "10,504", "10,504")),
    @("H", @("Tier 1 capital ratio (%)", "16.7", "16.3", "17.3")),
    @("I", @("Wages and salaries", "549", "535", "550")),
    @("J", @("Occupational pension provision and other benefits", "87", "101", "93"))
)

for ($ci = 0; $ci -lt $columns.Length; $ci++) {
    $colLetter = $columns[$ci][0]
    $values = $columns[$ci][1]
    for ($ri = 0; $ri -lt $values.Length; $ri++) {
        $cell = $ws.Range($colLetter + ($ri + 1))
        $cell.NumberFormat = "@"
        $cell.Value = $values[$ri]
        $cell.Style = "Normal"
    }
}
